$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.904.15"
Set-TextValue "E2" "  +0.46%  "
Set-TextValue "D3" "1.553.63"
Set-TextValue "E3" "  +0.17%  "
Set-TextValue "E4" "  +0.37%  "
Set-TextValue "D5" "206.50"
Set-TextValue "E5" "  +1.17%  "
Set-TextValue "D6" "0.485"
Set-TextValue "E6" "  +0.86%  "
Set-TextValue "E7" "  +0.36%  "
Set-TextValue "E8" "  +0.78%  "
Set-TextValue "D9" "21.47"
Set-TextValue "E9" "  +0.21%  "
Set-TextValue "D10" "0.0582"
Set-TextValue "E10" "  +0.06%  "
Set-TextValue "E11" "  +0.06%  "
Set-TextValue "D12" "1.776.61"
Set-TextValue "E12" "  +0.21%  "
Set-TextValue "D13" "1.565.24"
Set-TextValue "E13" "  +0.40%  "
Set-TextValue "D14" "3.70"
Set-TextValue "E14" "  +1.08%  "
Set-TextValue "E15" "  +0.87%  "
Set-TextValue "D16" "26.926.62"
Set-TextValue "E16" "  +0.59%  "
Set-TextValue "D17" "61.61"
Set-TextValue "D18" "214.00"
Set-TextValue "E18" "  +0.00%  "
Set-TextValue "D19" "0.0₃0687"
Set-TextValue "E19" "  +1.18%  "
Set-TextValue "D20" "7.23"
Set-TextValue "E20" "  -0.50%  "
Set-TextValue "E21" "  +0.34%  "
Set-TextValue "D22" "4.04"
Set-TextValue "E22" "  -0.87%  "
Set-TextValue "D23" "9.18"
Set-TextValue "E23" "  +1.26%  "
Set-TextValue "E24" "  -1.73%  "
Set-TextValue "D25" "153.09"
Set-TextValue "E25" "  +0.93%  "
Set-TextValue "D26" "6.65"
Set-TextValue "E26" "  +2.11%  "
Set-TextValue "E27" "  +0.15%  "
Set-TextValue "E28" "  +0.35%  "
Set-TextValue "E29" "  +1.06%  "
Set-TextValue "E30" "  -0.58%  "
Set-TextValue "E31" "  -0.65%  "
Set-TextValue "E32" "  +2.11%  "
Set-TextValue "D33" "1.367.02"
Set-TextValue "E33" "  +0.12%  "
Set-TextValue "E34" "  +1.69%  "
Set-TextValue "D35" "1.54"
Set-TextValue "E35" "  +3.06%  "
Set-TextValue "D36" "0.970"
Set-TextValue "E36" "  +5.54%  "
Set-TextValue "E37" "  +0.42%  "
Set-TextValue "E38" "  +0.91%  "
Set-TextValue "D39" "0.521"
Set-TextValue "E39" "  +0.19%  "
Set-TextValue "D40" "0.807"
Set-TextValue "E40" "  +0.72%  "
Set-TextValue "E41" "  +0.37%  "
Set-TextValue "D42" "0.987"
Set-TextValue "E42" "  -0.25%  "
Set-TextValue "E43" "  -0.48%  "
Set-TextValue "E44" "  +3.32%  "
Set-TextValue "D45" "63.54"
Set-TextValue "E45" "  +1.06%  "
Set-TextValue "D46" "1.73"
Set-TextValue "E46" "  -1.84%  "
Set-TextValue "D47" "1.689.35"
Set-TextValue "E47" "  +0.01%  "
Set-TextValue "D48" "86.10"
Set-TextValue "E48" "  +0.06%  "
Set-TextValue "E49" "  -0.60%  "
Set-TextValue "D50" "0.0955"
Set-TextValue "E50" "  +1.29%  "
Set-TextValue "E51" "  +0.50%  "
